# Generate Report for Handoff
# Marks file "b.md" as ready for handoff (instead of "handed back: in sync"),
# across the Overview / zh-cn / de-de sheets, and records the new handoff
# file names + timestamps for the zh-cn and de-de target languages.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet — row for b.md (row 3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-25-19 12:25:16"

# ---------------------------------------------------------------------
# zh-cn sheet — row for b.md (row 3)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 12:25:13"

# The engine's Hyperlinks collection only supports a sheet-wide Delete,
# so capture every existing hyperlink on this sheet (ref/address/display)
# before clearing, then replay them back with the updated display text
# for the one that changed (D3).
$zhCnLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/a.md"; Display = "a.md" },
    @{ Ref = "B2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/a.md"; Display = ".md" },
    @{ Ref = "D2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b61045e34e385beb0517d6bfcaaad0ad85c997b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/49beeeaf5f6a6495ffa8cbc77d9f345ff3d4242a/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ce07d17c5d0969a8acc09a7015f0775ada8c2d6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/b.md"; Display = "b.md" },
    @{ Ref = "B3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/b.md"; Display = ".md" },
    @{ Ref = "D3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b61045e34e385beb0517d6bfcaaad0ad85c997b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/49beeeaf5f6a6495ffa8cbc77d9f345ff3d4242a/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ce07d17c5d0969a8acc09a7015f0775ada8c2d6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" }
)

$wsZhCn.Hyperlinks.Delete()
foreach ($link in $zhCnLinks) {
    $target = $wsZhCn.Range($link.Ref)
    $wsZhCn.Hyperlinks.Add($target, $link.Address, [Type]::Missing, [Type]::Missing, $link.Display)
}

# ---------------------------------------------------------------------
# de-de sheet — row for b.md (row 3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 12:25:16"

$deDeLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/a.md"; Display = "a.md" },
    @{ Ref = "B2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/a.md"; Display = ".md" },
    @{ Ref = "D2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0552c4247855971843ad6b78e66d7171f7ea9bf8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0e9d9d6ac8065fc04f79e9668cc4a1815612fa48/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3b6ad859ae76171bdb9f572bed5faa39c5855a91/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/b.md"; Display = "b.md" },
    @{ Ref = "B3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/8108453b29e78b387155761a1a5e7cd6f710e24b/e2e/b.md"; Display = ".md" },
    @{ Ref = "D3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0552c4247855971843ad6b78e66d7171f7ea9bf8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0e9d9d6ac8065fc04f79e9668cc4a1815612fa48/e2e/a.md"; Display = "a.md" },
    @{ Ref = "G3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3b6ad859ae76171bdb9f572bed5faa39c5855a91/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" }
)

$wsDeDe.Hyperlinks.Delete()
foreach ($link in $deDeLinks) {
    $target = $wsDeDe.Range($link.Ref)
    $wsDeDe.Hyperlinks.Add($target, $link.Address, [Type]::Missing, [Type]::Missing, $link.Display)
}
